{"js": "async (context) => {\n  const body = context.document.body;\n  const paragraphs = body.paragraphs;\n  paragraphs.load(\"items\");\n  await context.sync();\n\n  // The very first paragraph in the document carries the bookmark\n  // \"_Hlk162277642\" and originally held the \"# CLL+1\" heading text.\n  // The edit pushes that heading further down and turns this first\n  // paragraph into the opening \"---\" of a new YAML front-matter block,\n  // followed by a description line, closing \"---\", a blank line, four\n  // import statements, another blank line and finally the \"# CLL+1\"\n  // heading (now as its own paragraph).\n  const firstPara = paragraphs.items[0];\n\n  // The original first paragraph (still bearing the bookmark) becomes\n  // the opening \"---\" delimiter of the new front-matter block. Replace\n  // the paragraph's range content (rather than using .clear(), which\n  // would also strip the bookmarkStart element) so the bookmark stays.\n  const firstRange = firstPara.getRange();\n  firstRange.insertText(\"---\", Word.InsertLocation.replace);\n\n  // Build the remaining new paragraphs in order, inserting each one\n  // directly after the previous one so the final order reads naturally\n  // from top to bottom.\n  const descPara = firstPara.insertParagraph(\"\", Word.InsertLocation.after);\n  descPara.insertText(\"description: History of the \", Word.InsertLocation.end);\n  descPara.insertText(\"CLL+1\", Word.InsertLocation.end);\n  descPara.insertText(\" step for Rubik's Cube.\", Word.InsertLocation.end);\n\n  const closeFrontMatter = descPara.insertParagraph(\"---\", Word.InsertLocation.after);\n\n  // Blank spacer paragraph.\n  const blank1 = closeFrontMatter.insertParagraph(\"\", Word.InsertLocation.after);\n\n  const importTwisty = blank1.insertParagraph(\"\", Word.InsertLocation.after);\n  importTwisty.insertText(\"import \", Word.InsertLocation.end);\n  importTwisty.insertText(\"TwistyPlayer\", Word.InsertLocation.end);\n  importTwisty.insertText(\" from \\\"@site/\", Word.InsertLocation.end);\n  importTwisty.insertText(\"src\", Word.InsertLocation.end);\n  importTwisty.insertText(\"/components/\", Word.InsertLocation.end);\n  importTwisty.insertText(\"TwistyPlayer\", Word.InsertLocation.end);\n  importTwisty.insertText(\"\\\";\", Word.InsertLocation.end);\n\n  const importExhibit = importTwisty.insertParagraph(\"\", Word.InsertLocation.after);\n  importExhibit.insertText(\"import Exhibit from \\\"@site/\", Word.InsertLocation.end);\n  importExhibit.insertText(\"src\", Word.InsertLocation.end);\n  importExhibit.insertText(\"/components/Exhibit\\\";\", Word.InsertLocation.end);\n\n  const importYouTube = importExhibit.insertParagraph(\"\", Word.InsertLocation.after);\n  importYouTube.insertText(\"import YouTube from \\\"@site/\", Word.InsertLocation.end);\n  importYouTube.insertText(\"src\", Word.InsertLocation.end);\n  importYouTube.insertText(\"/components/YouTube\\\";\", Word.InsertLocation.end);\n\n  const importCollage = importYouTube.insertParagraph(\"\", Word.InsertLocation.after);\n  importCollage.insertText(\"import \", Word.InsertLocation.end);\n  importCollage.insertText(\"ImageCollage\", Word.InsertLocation.end);\n  importCollage.insertText(\" from '@site/\", Word.InsertLocation.end);\n  importCollage.insertText(\"src\", Word.InsertLocation.end);\n  importCollage.insertText(\"/components/\", Word.InsertLocation.end);\n  importCollage.insertText(\"ImageCollage\", Word.InsertLocation.end);\n  importCollage.insertText(\"';\", Word.InsertLocation.end);\n\n  // Blank spacer paragraph.\n  const blank2 = importCollage.insertParagraph(\"\", Word.InsertLocation.after);\n\n  const headingPara = blank2.insertParagraph(\"# CLL+1\", Word.InsertLocation.after);\n\n  await context.sync();\n\n  // Re-load paragraphs to locate the blank paragraph that originally sat\n  // right before \"## Description\" -- the new <Exhibit ... /> block is\n  // inserted right after it.\n  const paragraphs2 = body.paragraphs;\n  paragraphs2.load(\"items/text\");\n  await context.sync();\n\n  let descriptionHeadingPara = null;\n  for (let i = 0; i < paragraphs2.items.length; i++) {\n    if (paragraphs2.items[i].text === \"## Description\") {\n      descriptionHeadingPara = paragraphs2.items[i];\n      break;\n    }\n  }\n\n  if (descriptionHeadingPara) {\n    const exhibitOpen = descriptionHeadingPara.insertParagraph(\"<Exhibit\", Word.InsertLocation.before);\n\n    const stickeringPara = descriptionHeadingPara.insertParagraph(\"\", Word.InsertLocation.before);\n    stickeringPara.insertText(\"  stickering\", Word.InsertLocation.end);\n    stickeringPara.insertText(\"={\", Word.InsertLocation.end);\n    stickeringPara.insertText(\"{\", Word.InsertLocation.end);\n\n    const solvedPara = descriptionHeadingPara.insertParagraph(\"\", Word.InsertLocation.before);\n    solvedPara.insertText(\"    solved: \\\"U D F B L R D\", Word.InsertLocation.end);\n    solvedPara.insertText(\"F DL DB DR DFR DFL DBL DBR FR FL BL BR UFR UFL UBL UBR UB\", Word.InsertLocation.end);\n    solvedPara.insertText(\"\\\"}}\", Word.InsertLocation.end);\n\n    const closeTag = descriptionHeadingPara.insertParagraph(\"/>\", Word.InsertLocation.before);\n\n    // Blank spacer paragraph between the Exhibit block and \"## Description\".\n    descriptionHeadingPara.insertParagraph(\"\", Word.InsertLocation.before);\n  }\n\n  await context.sync();\n};\n", "ps1": "$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# 1) Front-matter + import block.\n#\n# The very first paragraph in the document carries the bookmark\n# \"_Hlk162277642\" and originally held the \"# CLL+1\" heading text. The\n# edit turns this first paragraph into the opening \"---\" of a new YAML\n# front-matter block (keeping the bookmark in place) and adds, right\n# after it: a description line, the closing \"---\", a blank line, four\n# import statements, another blank line and finally the \"# CLL+1\"\n# heading as its own paragraph.\n# ---------------------------------------------------------------------\n\n$firstPara = $d.Paragraphs.Item(1)\n\n# Replacing the Range text (instead of e.g. typing over a selection)\n# keeps the existing bookmarkStart/bookmarkEnd pair attached to this\n# paragraph.\n$firstPara.Range.Text = \"---\"\n\n# Insert 9 blank paragraphs directly after the (now \"---\") first\n# paragraph; fill them in afterwards by final position so the ordering\n# reads top-to-bottom in call order.\nfor ($i = 0; $i -lt 9; $i++) {\n  $firstPara.Range.InsertParagraphAfter()\n}\n\n$d.Paragraphs.Item(2).Range.Text = \"description: History of the CLL+1 step for Rubik's Cube.\"\n$d.Paragraphs.Item(3).Range.Text = \"---\"\n# Paragraph 4 stays blank.\n$d.Paragraphs.Item(5).Range.Text = 'import TwistyPlayer from \"@site/src/components/TwistyPlayer\";'\n$d.Paragraphs.Item(6).Range.Text = 'import Exhibit from \"@site/src/components/Exhibit\";'\n$d.Paragraphs.Item(7).Range.Text = 'import YouTube from \"@site/src/components/YouTube\";'\n$d.Paragraphs.Item(8).Range.Text = \"import ImageCollage from '@site/src/components/ImageCollage';\"\n# Paragraph 9 stays blank.\n$d.Paragraphs.Item(10).Range.Text = \"# CLL+1\"\n\n# ---------------------------------------------------------------------\n# 2) <Exhibit ... /> component block, inserted right before\n#    \"## Description\" (after the blank paragraph that already precedes\n#    it).\n# ---------------------------------------------------------------------\n\n$findRange = $d.Content\n$findRange.Find.ClearFormatting()\n$findRange.Find.Execute(\"## Description\") | Out-Null\n$descPara = $findRange.Paragraphs.Item(1)\n\nfor ($i = 0; $i -lt 5; $i++) {\n  $descPara.Range.InsertParagraphBefore()\n}\n\n$descIndex = $descPara.Range.Paragraphs.Item(1).Index\n$exhibitStart = $descIndex - 5\n\n$d.Paragraphs.Item($exhibitStart).Range.Text = \"<Exhibit\"\n$d.Paragraphs.Item($exhibitStart + 1).Range.Text = \"  stickering={{\"\n$d.Paragraphs.Item($exhibitStart + 2).Range.Text = '    solved: \"U D F B L R DF DL DB DR DFR DFL DBL DBR FR FL BL BR UFR UFL UBL UBR UB\"}}'\n$d.Paragraphs.Item($exhibitStart + 3).Range.Text = \"/>\"\n# $exhibitStart + 4 stays blank.\n\nWrite-Output \"done\"\n"}
